$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04978849882897407
$ws.Range("D2").Value = 0.01742262733052513
$ws.Range("E2").Value = 0.1258385180787442
$ws.Range("F2").Value = 0.4881132597339999
$ws.Range("G2").Value = 0.3330331740223045
$ws.Range("H2").Value = 0.4851437576033106
$ws.Range("K2").Value = 0.7882615431310853
$ws.Range("M2").Value = 0.3142111624823514
$ws.Range("N2").Value = 1.003008170214791
$ws.Range("O2").Value = 1.569274099121031
$ws.Range("C3").Value = 0.04419155087620652
$ws.Range("D3").Value = 0.0155700706879216
$ws.Range("E3").Value = 0.1197778099722839
$ws.Range("F3").Value = 0.4827594287832468
$ws.Range("G3").Value = 0.328925592314647
$ws.Range("H3").Value = 0.4871128233655213
$ws.Range("K3").Value = 0.6877651373429501
$ws.Range("M3").Value = 0.2795813537323824
$ws.Range("N3").Value = 1.006940546819074
$ws.Range("O3").Value = 1.564406486079974
$ws.Range("C4").Value = 0.04076807626913137
$ws.Range("D4").Value = 0.01442565598295431
$ws.Range("E4").Value = 0.1161755229688168
$ws.Range("F4").Value = 0.479850183437101
$ws.Range("G4").Value = 0.326710464224206
$ws.Range("H4").Value = 0.4885973284187486
$ws.Range("K4").Value = 0.6258494318049088
$ws.Range("M4").Value = 0.2583655859229026
$ws.Range("N4").Value = 1.009755151966765
$ws.Range("O4").Value = 1.562635162836003
$ws.Range("C5").Value = 0.03937627889396822
$ws.Range("D5").Value = 0.01395758438815164
$ws.Range("E5").Value = 0.1147372319677231
$ws.Range("F5").Value = 0.4787595673505791
$ws.Range("G5").Value = 0.3258847335387998
$ws.Range("H5").Value = 0.4892715279916757
$ws.Range("K5").Value = 0.6005669088142724
$ws.Range("M5").Value = 0.2497319526691513
$ws.Range("N5").Value = 1.011002908307816
$ws.Range("O5").Value = 1.562218916976605
$ws.Range("C6").Value = 0.03914537135071328
$ws.Range("D6").Value = 0.01387975889493021
$ws.Range("E6").Value = 0.114500189355077
$ws.Range("F6").Value = 0.4785842007141881
$ws.Range("G6").Value = 0.3257522625133333
$ws.Range("H6").Value = 0.4893876602947742
$ws.Range("K6").Value = 0.5963657085929697
$ws.Range("M6").Value = 0.2482990708627568
$ws.Range("N6").Value = 1.011216189405609
$ws.Range("O6").Value = 1.562168240154918
$ws.Range("C7").Value = 0.04074929264342586
$ws.Range("D7").Value = 0.01441935030059938
$ws.Range("E7").Value = 0.1161560058786435
$ws.Range("F7").Value = 0.4798350908371845
$ws.Range("G7").Value = 0.3266990168632944
$ws.Range("H7").Value = 0.4886061405434958
$ws.Range("K7").Value = 0.6255086687793323
$ws.Range("M7").Value = 0.2582491011119856
$ws.Range("N7").Value = 1.009771571354094
$ws.Range("O7").Value = 1.562628312632242
$ws.Range("C8").Value = 0.04785596126501446
$ws.Range("D8").Value = 0.01678532381384912
$ws.Range("E8").Value = 0.1237239285507812
$ws.Range("F8").Value = 0.4861887157388765
$ws.Range("G8").Value = 0.3315530322742148
$ws.Range("H8").Value = 0.4857654989000508
$ws.Range("K8").Value = 0.7536549796150211
$ws.Range("M8").Value = 0.3022610993565777
$ws.Range("N8").Value = 1.004281110727874
$ws.Range("O8").Value = 1.567342738258588
$ws.Range("C9").Value = 0.06189617685159021
$ws.Range("D9").Value = 0.02136878334593462
$ws.Range("E9").Value = 0.1395213402769429
$ws.Range("F9").Value = 0.5016550989782402
$ws.Range("G9").Value = 0.3435186767894862
$ws.Range("H9").Value = 0.482382350686521
$ws.Range("K9").Value = 1.003222368778722
$ws.Range("M9").Value = 0.3889408711203259
$ws.Range("N9").Value = 0.9966822600424337
$ws.Range("O9").Value = 1.586274988567453
$ws.Range("C10").Value = 0.07227644189245552
$ws.Range("D10").Value = 0.02470072906987042
$ws.Range("E10").Value = 0.1517293986921473
$ws.Range("F10").Value = 0.5148637661094853
$ws.Range("G10").Value = 0.3538189202284769
$ws.Range("H10").Value = 0.4812328279991362
$ws.Range("K10").Value = 1.185467595210071
$ws.Range("M10").Value = 0.4528565501416182
$ws.Range("N10").Value = 0.9930226314058359
$ws.Range("O10").Value = 1.606132864313366
$ws.Range("C11").Value = 0.07701316080856202
$ws.Range("D11").Value = 0.02620855601842464
$ws.Range("E11").Value = 0.1574178577494152
$ws.Range("F11").Value = 0.5212762455103501
$ws.Range("G11").Value = 0.3588362583106175
$ws.Range("H11").Value = 0.4810006056422509
$ws.Range("K11").Value = 1.268122799173852
$ws.Range("M11").Value = 0.481985318237534
$ws.Range("N11").Value = 0.9917738599422563
$ws.Range("O11").Value = 1.616467793276485
$ws.Range("C12").Value = 0.07880894957970952
$ws.Range("D12").Value = 0.02677836814428503
$ws.Range("E12").Value = 0.1595916350107203
$ws.Range("F12").Value = 0.5237627280996975
$ws.Range("G12").Value = 0.3607841664688038
$ws.Range("H12").Value = 0.4809545117952183
$ws.Range("K12").Value = 1.299385007784906
$ws.Range("M12").Value = 0.4930232467596767
$ws.Range("N12").Value = 0.9913606766919969
$ws.Range("O12").Value = 1.620569185373114
$ws.Range("C13").Value = 0.07842210136438155
$ws.Range("D13").Value = 0.02665570145543228
$ws.Range("E13").Value = 0.1591225934657743
$ws.Range("F13").Value = 0.5232246280284443
$ws.Range("G13").Value = 0.3603625129534862
$ws.Range("H13").Value = 0.4809625772627157
$ws.Range("K13").Value = 1.29265382781233
$ws.Range("M13").Value = 0.4906457004205151
$ws.Range("N13").Value = 0.9914470098437249
$ws.Range("O13").Value = 1.619677515731127
$ws.Range("C14").Value = 0.07716085951383889
$ws.Range("D14").Value = 0.02625545843464749
$ws.Range("E14").Value = 0.1575962996925568
$ws.Range("F14").Value = 0.5214796423368142
$ws.Range("G14").Value = 0.3589955514354131
$ws.Range("H14").Value = 0.4809959746273051
$ws.Range("K14").Value = 1.270695521681262
$ws.Range("M14").Value = 0.4828932659091691
$ws.Range("N14").Value = 0.9917386713253649
$ws.Range("O14").Value = 1.616801450323635
$ws.Range("C15").Value = 0.07638858496908085
$ws.Range("D15").Value = 0.02601014474193875
$ws.Range("E15").Value = 0.1566639719460099
$ws.Range("F15").Value = 0.5204183740959394
$ws.Range("G15").Value = 0.3581645001494991
$ws.Range("H15").Value = 0.4810218819886103
$ws.Range("K15").Value = 1.257240487889135
$ws.Range("M15").Value = 0.478145647570031
$ws.Range("N15").Value = 0.9919250934845394
$ws.Range("O15").Value = 1.615064250125982
$ws.Range("C16").Value = 0.07196718034740002
$ws.Range("D16").Value = 0.02460202728802585
$ws.Range("E16").Value = 0.1513603836230146
$ws.Range("F16").Value = 0.5144528337568417
$ws.Range("G16").Value = 0.3534977217451853
$ws.Range("H16").Value = 0.4812538574534528
$ws.Range("K16").Value = 1.180060725309602
$ws.Range("M16").Value = 0.4509539740416244
$ws.Range("N16").Value = 0.9931126021730847
$ws.Range("O16").Value = 1.605483695230589
$ws.Range("C17").Value = 0.06925854199054982
$ws.Range("D17").Value = 0.02373614706905869
$ws.Range("E17").Value = 0.1481415642195572
$ws.Range("F17").Value = 0.510896696468599
$ws.Range("G17").Value = 0.3507199420443072
$ws.Range("H17").Value = 0.4814706517968546
$ws.Range("K17").Value = 1.132648432649376
$ws.Range("M17").Value = 0.4342862800732092
$ws.Range("N17").Value = 0.9939475619201517
$ws.Range("O17").Value = 1.599940123085901
$ws.Range("C18").Value = 0.06770198796246518
$ws.Range("D18").Value = 0.02323737489859923
$ws.Range("E18").Value = 0.1463028702815805
$ws.Range("F18").Value = 0.5088892980744006
$ws.Range("G18").Value = 0.349153443343809
$ws.Range("H18").Value = 0.4816227058808664
$ws.Range("K18").Value = 1.105354794669381
$ws.Range("M18").Value = 0.4247044804933466
$ws.Range("N18").Value = 0.9944669726631048
$ws.Range("O18").Value = 1.596874071299283
$ws.Range("C19").Value = 0.0671752035273272
$ws.Range("D19").Value = 0.0230683730330199
$ws.Range("E19").Value = 0.1456824905297793
$ws.Range("F19").Value = 0.5082161489845163
$ws.Range("G19").Value = 0.3486284058760134
$ws.Range("H19").Value = 0.4816788865829551
$ws.Range("K19").Value = 1.096109681068754
$ws.Range("M19").Value = 0.4214611142081708
$ws.Range("N19").Value = 0.9946495664977419
$ws.Range("O19").Value = 1.595856973049109
$ws.Range("C20").Value = 0.06954673813007162
$ws.Range("D20").Value = 0.02382839841197182
$ws.Range("E20").Value = 0.1484828981099398
$ws.Range("F20").Value = 0.5112713196319802
$ws.Range("G20").Value = 0.3510124098260974
$ws.Range("H20").Value = 0.4814447417419103
$ws.Range("K20").Value = 1.137697980263624
$ws.Range("M20").Value = 0.4360600664170846
$ws.Range("N20").Value = 0.9938546266083961
$ws.Range("O20").Value = 1.600517566653537
$ws.Range("C21").Value = 0.07753126005673039
$ws.Range("D21").Value = 0.0263730515181777
$ws.Range("E21").Value = 0.1580440728314514
$ws.Range("F21").Value = 0.5219906057750023
$ws.Range("G21").Value = 0.3593957577859896
$ws.Range("H21").Value = 0.4809850290568818
$ws.Range("K21").Value = 1.277146240429545
$ws.Range("M21").Value = 0.4851701416323237
$ws.Range("N21").Value = 0.99165138400312
$ws.Range("O21").Value = 1.617641119203142
$ws.Range("C22").Value = 0.08276181689713269
$ws.Range("D22").Value = 0.0280292984230428
$ws.Range("E22").Value = 0.1644076987794065
$ws.Range("F22").Value = 0.5293356743857345
$ws.Range("G22").Value = 0.3651543721975656
$ws.Range("H22").Value = 0.4809284965974143
$ws.Range("K22").Value = 1.368064236027124
$ws.Range("M22").Value = 0.5173100369327841
$ws.Range("N22").Value = 0.9905593584130514
$ws.Range("O22").Value = 1.629927158225456
$ws.Range("C23").Value = 0.07996905829364209
$ws.Range("D23").Value = 0.02714596466051233
$ws.Range("E23").Value = 0.1610007124011759
$ws.Range("F23").Value = 0.5253843687306841
$ws.Range("G23").Value = 0.3620552238124475
$ws.Range("H23").Value = 0.480936337051034
$ws.Range("K23").Value = 1.319560238449469
$ws.Range("M23").Value = 0.5001524295051354
$ws.Range("N23").Value = 0.9911103973218616
$ws.Range("O23").Value = 1.623269486509599
$ws.Range("C24").Value = 0.06941644254150958
$ws.Range("D24").Value = 0.02378669458666138
$ws.Range("E24").Value = 0.1483285441746887
$ws.Range("F24").Value = 0.5111018370334151
$ws.Range("G24").Value = 0.3508800902059903
$ws.Range("H24").Value = 0.4814563702797017
$ws.Range("K24").Value = 1.135415190796323
$ws.Range("M24").Value = 0.4352581354859097
$ws.Range("N24").Value = 0.9938965199806802
$ws.Range("O24").Value = 1.600256127486148
$ws.Range("C25").Value = 0.05808663120554058
$ws.Range("D25").Value = 0.0201349835719995
$ws.Range("E25").Value = 0.1351433910557773
$ws.Range("F25").Value = 0.4971478320345923
$ws.Range("G25").Value = 0.3400178921582011
$ws.Range("H25").Value = 0.4830631421126128
$ws.Range("K25").Value = 0.9358988005172364
$ws.Range("M25").Value = 0.3654513113842697
$ws.Range("N25").Value = 0.9966822600424337
$ws.Range("O25").Value = 1.586274988567453
